$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.830478549003601
$ws.Range("B1").Value = 1.927390098571777
$ws.Range("C1").Value = 2.105088710784912
$ws.Range("D1").Value = 2.975113153457642
$ws.Range("E1").Value = 4.216721534729004
